$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(10002, 110021, 10021),
    @(10003, 110022, 10022),
    @(10004, 110023, 10023),
    @(10005, 110024, 10024),
    @(10006, 110025, 10025),
    @(10007, 110026, 10026),
    @(10008, 110027, 10027),
    @(10009, 110028, 10028),
    @(10010, 110029, 10029)
)

$row = 22
foreach ($d in $data) {
    $ws.Cells.Item($row, 1).Value = $d[0]
    $ws.Cells.Item($row, 2).Value = $d[1]
    $ws.Cells.Item($row, 3).Value = $d[2]
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin"
    $ws.Cells.Item($row, 7).Value = "now()"
    $ws.Cells.Item($row, 8).Value = "now()"
    $row++
}

$ws.Range("F14").Select()

# Page setup: printed in portrait orientation at 300 dpi (as set via Page Setup dialog)
$xlPortrait = 1
$ws.PageSetup.Orientation = $xlPortrait
$ws.PageSetup.PrintQuality = 300
